$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "add report offset range"
#
# Populate the dcd_error_min (%) / dcd_error_max (%) columns (V, W) with the
# measured report offset range instead of the "error" / "n/a" placeholder
# text that filled those cells before.

# Rows 2-10: per-corner data rows. They previously carried the yellow
# "error" placeholder style (s="2") on V/W; clear that formatting now that
# real numbers are present so the cells fall back to the default style.
$ws.Range("V2:W10").ClearFormats()

$ws.Range("V2").Value = -6.66624
$ws.Range("W2").Value = 6.65638

$ws.Range("V3").Value = -7.9208
$ws.Range("W3").Value = 7.85949

$ws.Range("V4").Value = -6.04503
$ws.Range("W4").Value = 6.07789

$ws.Range("V5").Value = -7.35354
$ws.Range("W5").Value = 7.36635

$ws.Range("V6").Value = -5.47758
$ws.Range("W6").Value = 5.47207

$ws.Range("V7").Value = -7.53332
$ws.Range("W7").Value = 7.57574

$ws.Range("V8").Value = -5.61293
$ws.Range("W8").Value = 5.61112

$ws.Range("V9").Value = -7.40006
$ws.Range("W9").Value = 7.33337

$ws.Range("V10").Value = -5.51945
$ws.Range("W10").Value = 5.56523

# Rows 12, 13, 18, 19, 20: summary block (Min / Max / Internal bounds / Mean
# / Std. Dev. / Rel. Std. Dev.). These keep their bold custom-format style
# (s="1") - only the "n/a" placeholder values in V/W are replaced with the
# computed numbers.
$ws.Range("V12").Value = -7.9208
$ws.Range("W12").Value = 5.47207

$ws.Range("V13").Value = -5.47758
$ws.Range("W13").Value = 7.85949

$ws.Range("V18").Value = -6.614328
$ws.Range("W18").Value = 6.613071

$ws.Range("V19").Value = 0.915045
$ws.Range("W19").Value = 0.89806

$ws.Range("V20").Value = 13.834285
$ws.Range("W20").Value = 13.580075
